$d = $word.ActiveDocument

# Locate the "GIS & Geospatial Analysis Consulting" paragraph under the
# Siege Analytics / PARTNER heading so we insert the new bullet points
# right after it (and before the existing "Lead comprehensive research..." bullet).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "GIS & Geospatial Analysis Consulting`r") {
        $targetIndex = $i
        break
    }
}

$bullets = @(
    "• Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels",
    "• Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide",
    "• Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis"
)

$insertAfterIndex = $targetIndex
foreach ($bullet in $bullets) {
    $anchor = $d.Paragraphs.Item($insertAfterIndex)
    $rng = $anchor.Range
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()
    $insertAfterIndex = $insertAfterIndex + 1
    $newPara = $d.Paragraphs.Item($insertAfterIndex)
    $newPara.Range.Text = $bullet
}
